$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 106250
$ws.Range("J3").Value = 106250
$ws.Range("L3").Value = 106250
$ws.Range("N3").Value = -106478
$ws.Range("H12").Value = 282.66666
$ws.Range("I12").Value = 282.66666
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 282.66666
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -112.66666
$ws.Range("N12").ClearContents()
$ws.Range("H32").Value = 3486.6365
$ws.Range("I32").Value = 3192
$ws.Range("J32").Value = 3732.1667
$ws.Range("K32").Value = 3192
$ws.Range("L32").Value = 3732.1667
$ws.Range("M32").Value = -2866
$ws.Range("N32").Value = -4384.1667
$ws.Range("H34").Value = 11428.25
$ws.Range("I34").Value = 7269
$ws.Range("J34").Value = 32224.5
$ws.Range("K34").Value = 7269
$ws.Range("L34").Value = 32224.5
$ws.Range("M34").Value = -7066
$ws.Range("N34").Value = -32630.5
$ws.Range("H36").Value = 11428.25
$ws.Range("I36").Value = 7269
$ws.Range("J36").Value = 32224.5
$ws.Range("K36").Value = 7269
$ws.Range("L36").Value = 32224.5
$ws.Range("M36").Value = -6554
$ws.Range("N36").Value = -33654.5
$ws.Range("H40").Value = 3499.5
$ws.Range("I40").Value = 3499.5
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 3499.5
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -3324.5
$ws.Range("N40").ClearContents()
$ws.Range("H80").Value = 67983.44500000001
$ws.Range("I80").Value = 200700
$ws.Range("J80").Value = 1625.1666
$ws.Range("K80").Value = 602100
$ws.Range("L80").Value = 4875.4998
$ws.Range("M80").Value = -601102
$ws.Range("N80").Value = -6871.4998
$ws.Range("H83").Value = 67983.44500000001
$ws.Range("I83").Value = 200700
$ws.Range("J83").Value = 1625.1666
$ws.Range("K83").Value = 1806300
$ws.Range("L83").Value = 14626.4994
$ws.Range("M83").Value = -1801308
$ws.Range("N83").Value = -24610.4994
$ws.Range("H86").Value = 61478.41
$ws.Range("I86").Value = 79925.62
$ws.Range("J86").Value = 1525
$ws.Range("K86").Value = 79925.62
$ws.Range("L86").Value = 1525
$ws.Range("M86").Value = -78802.62
$ws.Range("N86").Value = -3771
$ws.Range("H89").Value = 61478.41
$ws.Range("I89").Value = 79925.62
$ws.Range("J89").Value = 1525
$ws.Range("K89").Value = 399628.1
$ws.Range("L89").Value = 7625
$ws.Range("M89").Value = -394012.1
$ws.Range("N89").Value = -18857
$ws.Range("H102").Value = 106250
$ws.Range("J102").Value = 106250
$ws.Range("L102").Value = 106250
$ws.Range("N102").Value = -112740
$ws.Range("H105").Value = 24632.666
$ws.Range("J105").Value = 24632.666
$ws.Range("L105").Value = 24632.666
$ws.Range("N105").Value = -31620.666
$ws.Range("H112").Value = 3863.818
$ws.Range("J112").Value = 3863.818
$ws.Range("L112").Value = 11591.454
$ws.Range("N112").Value = -13807.454
$ws.Range("H125").Value = 1114.7142
$ws.Range("J125").Value = 1243.8334
$ws.Range("L125").Value = 11194.5006
$ws.Range("N125").Value = -16114.5006
$ws.Range("H129").Value = 7434.9
$ws.Range("I129").Value = 2527.1428
$ws.Range("K129").Value = 7581.428400000001
$ws.Range("M129").Value = -2581.428400000001
$ws.Range("H132").Value = 4901.3955
$ws.Range("I132").Value = 2783.9473
$ws.Range("K132").Value = 8351.841899999999
$ws.Range("M132").Value = -5821.841899999999
$ws.Range("H135").Value = 2832
$ws.Range("I135").Value = 1670.75
$ws.Range("K135").Value = 15036.75
$ws.Range("M135").Value = -12501.75
$ws.Range("H138").Value = 3514.3396
$ws.Range("I138").Value = 3107.3333
$ws.Range("J138").Value = 3597.5908
$ws.Range("K138").Value = 9321.999899999999
$ws.Range("L138").Value = 10792.7724
$ws.Range("M138").Value = -4181.999899999999
$ws.Range("N138").Value = -21072.7724

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7079.88
$ws.Range("I32").Value = 4476.1904
$ws.Range("K32").Value = 4476.1904
$ws.Range("M32").Value = -4189.1904
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H88").Value = 1945.8667
$ws.Range("I88").Value = 1460.625
$ws.Range("J88").Value = 2500.4285
$ws.Range("K88").Value = 1460.625
$ws.Range("L88").Value = 2500.4285
$ws.Range("M88").Value = -1054.625
$ws.Range("N88").Value = -3312.4285
$ws.Range("H91").Value = 1945.8667
$ws.Range("I91").Value = 1460.625
$ws.Range("J91").Value = 2500.4285
$ws.Range("K91").Value = 1460.625
$ws.Range("L91").Value = 2500.4285
$ws.Range("M91").Value = -56.625
$ws.Range("N91").Value = -5308.4285
$ws.Range("H133").Value = 58000
$ws.Range("J133").Value = 58000
$ws.Range("L133").Value = 58000
$ws.Range("N133").Value = -63060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2726.75
$ws.Range("I20").Value = 1953.75
$ws.Range("J20").Value = 3499.75
$ws.Range("K20").Value = 1953.75
$ws.Range("L20").Value = 3499.75
$ws.Range("M20").Value = -1706.75
$ws.Range("N20").Value = -3993.75
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H48").Value = 0
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()
$ws.Range("H76").Value = 28602.666
$ws.Range("J76").Value = 28602.666
$ws.Range("L76").Value = 28602.666
$ws.Range("N76").Value = -29232.666
$ws.Range("H79").Value = 28602.666
$ws.Range("J79").Value = 28602.666
$ws.Range("L79").Value = 28602.666
$ws.Range("N79").Value = -30786.666
$ws.Range("H94").Value = 7193.826
$ws.Range("I94").Value = 8292.833000000001
$ws.Range("J94").Value = 3237.4
$ws.Range("K94").Value = 8292.833000000001
$ws.Range("L94").Value = 3237.4
$ws.Range("M94").Value = -7841.833000000001
$ws.Range("N94").Value = -4139.4
$ws.Range("H105").Value = 4529.4375
$ws.Range("I105").Value = 4410.0415
$ws.Range("K105").Value = 4410.0415
$ws.Range("M105").Value = -2663.0415
$ws.Range("H134").Value = 3607.3333
$ws.Range("I134").Value = 3201.8572
$ws.Range("J134").Value = 3962.125
$ws.Range("K134").Value = 9605.571599999999
$ws.Range("L134").Value = 11886.375
$ws.Range("M134").Value = -7070.571599999999
$ws.Range("N134").Value = -16956.375

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4117.3335
$ws.Range("I31").Value = 2518.0557
$ws.Range("J31").Value = 8915.166999999999
$ws.Range("K31").Value = 2518.0557
$ws.Range("L31").Value = 8915.166999999999
$ws.Range("M31").Value = -2223.0557
$ws.Range("N31").Value = -9505.166999999999
$ws.Range("H34").Value = 4117.3335
$ws.Range("I34").Value = 2518.0557
$ws.Range("J34").Value = 8915.166999999999
$ws.Range("K34").Value = 2518.0557
$ws.Range("L34").Value = 8915.166999999999
$ws.Range("M34").Value = -2316.0557
$ws.Range("N34").Value = -9319.166999999999
$ws.Range("H58").Value = 145983.58
$ws.Range("I58").Value = 169564.17
$ws.Range("K58").Value = 169564.17
$ws.Range("M58").Value = -169361.17
$ws.Range("H95").Value = 15966.667
$ws.Range("J95").Value = 15966.667
$ws.Range("L95").Value = 15966.667
$ws.Range("N95").Value = -21458.667
$ws.Range("H99").Value = 2648.9546
$ws.Range("I99").Value = 1995.0834
$ws.Range("J99").Value = 3433.6
$ws.Range("K99").Value = 1995.0834
$ws.Range("L99").Value = 3433.6
$ws.Range("M99").Value = -497.0834
$ws.Range("N99").Value = -6429.6
$ws.Range("H126").Value = 2648.9546
$ws.Range("I126").Value = 1995.0834
$ws.Range("J126").Value = 3433.6
$ws.Range("K126").Value = 5985.2502
$ws.Range("L126").Value = 10300.8
$ws.Range("M126").Value = -3515.2502
$ws.Range("N126").Value = -15240.8
$ws.Range("H134").Value = 37320.137
$ws.Range("I134").Value = 42887.04
$ws.Range("K134").Value = 128661.12
$ws.Range("M134").Value = -126126.12
$ws.Range("H136").Value = 145983.58
$ws.Range("I136").Value = 169564.17
$ws.Range("K136").Value = 508692.51
$ws.Range("M136").Value = -506142.51
$ws.Range("H141").Value = 299927.2
$ws.Range("I141").Value = 45777.6
$ws.Range("J141").Value = 511718.5
$ws.Range("K141").Value = 45777.6
$ws.Range("L141").Value = 511718.5
$ws.Range("M141").Value = -40597.6
$ws.Range("N141").Value = -522078.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 91.47369
$ws.Range("I2").Value = 39
$ws.Range("K2").Value = 234
$ws.Range("M2").Value = -121
$ws.Range("H5").Value = 935.0909
$ws.Range("I5").Value = 682.5
$ws.Range("K5").Value = 2047.5
$ws.Range("M5").Value = -1935.5
$ws.Range("H6").Value = 25.75
$ws.Range("I6").Value = 25.75
$ws.Range("K6").Value = 77.25
$ws.Range("M6").Value = 35.75
$ws.Range("H17").Value = 649.75
$ws.Range("I17").Value = 299.5
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 898.5
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -729.5
$ws.Range("N17").Value = -3338
$ws.Range("H23").Value = 71428650
$ws.Range("I23").Value = 79.40000000000001
$ws.Range("K23").Value = 238.2
$ws.Range("M23").Value = -3.200000000000017
$ws.Range("H34").Value = 2166.889
$ws.Range("J34").Value = 2400
$ws.Range("L34").Value = 7200
$ws.Range("N34").Value = -7368
$ws.Range("H39").Value = 2950
$ws.Range("I39").Value = 150
$ws.Range("J39").Value = 4350
$ws.Range("K39").Value = 450
$ws.Range("L39").Value = 13050
$ws.Range("M39").Value = -156
$ws.Range("N39").Value = -13638
$ws.Range("H40").Value = 178.66667
$ws.Range("I40").Value = 158.28572
$ws.Range("K40").Value = 633.14288
$ws.Range("M40").Value = -564.14288
$ws.Range("H55").Value = 7799.222
$ws.Range("H68").Value = 552.9474
$ws.Range("I68").Value = 575.1111
$ws.Range("K68").Value = 1725.3333
$ws.Range("M68").Value = -914.3332999999998
$ws.Range("H71").Value = 552.9474
$ws.Range("I71").Value = 575.1111
$ws.Range("K71").Value = 5175.9999
$ws.Range("M71").Value = -1119.9999
$ws.Range("H93").Value = 8644.833000000001
$ws.Range("I93").Value = 623
$ws.Range("K93").Value = 1869
$ws.Range("M93").Value = 3
$ws.Range("H120").Value = 19994
$ws.Range("I120").Value = 19994
$ws.Range("K120").Value = 59982
$ws.Range("M120").Value = -55144
$ws.Range("H135").Value = 935.0909
$ws.Range("I135").Value = 682.5
$ws.Range("K135").Value = 6142.5
$ws.Range("M135").Value = -3607.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1928.6
$ws.Range("I102").Value = 1800.7142
$ws.Range("J102").Value = 2440.1428
$ws.Range("K102").Value = 1800.7142
$ws.Range("L102").Value = 2440.1428
$ws.Range("M102").Value = -178.7141999999999
$ws.Range("N102").Value = -5684.1428
$ws.Range("H113").Value = 3229.5
$ws.Range("J113").Value = 3395.6
$ws.Range("L113").Value = 3395.6
$ws.Range("N113").Value = -7735.6
$ws.Range("H122").Value = 2220.5557
$ws.Range("I122").Value = 959.4
$ws.Range("J122").Value = 3797
$ws.Range("K122").Value = 2878.2
$ws.Range("L122").Value = 11391
$ws.Range("M122").Value = -428.1999999999998
$ws.Range("N122").Value = -16291
$ws.Range("H126").Value = 5175.7144
$ws.Range("I126").Value = 5057.7617
$ws.Range("K126").Value = 15173.2851
$ws.Range("M126").Value = -12703.2851
$ws.Range("H132").Value = 95596.91
$ws.Range("I132").Value = 113729.555
$ws.Range("J132").Value = 14000
$ws.Range("K132").Value = 341188.665
$ws.Range("L132").Value = 42000
$ws.Range("M132").Value = -338658.665
$ws.Range("N132").Value = -47060

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1162.6666
$ws.Range("I16").Value = 1162.6666
$ws.Range("K16").Value = 1162.6666
$ws.Range("M16").Value = -992.6666
$ws.Range("H46").Value = 19108.895
$ws.Range("I46").Value = 25267.154
$ws.Range("J46").Value = 5766
$ws.Range("K46").Value = 25267.154
$ws.Range("L46").Value = 5766
$ws.Range("M46").Value = -25079.154
$ws.Range("N46").Value = -6142
$ws.Range("H68").Value = 6742
$ws.Range("J68").Value = 7989.6665
$ws.Range("L68").Value = 7989.6665
$ws.Range("N68").Value = -9487.666499999999
$ws.Range("H71").Value = 6742
$ws.Range("J71").Value = 7989.6665
$ws.Range("L71").Value = 39948.3325
$ws.Range("N71").Value = -47436.3325
$ws.Range("H132").Value = 53274.375
$ws.Range("I132").Value = 73540.3
$ws.Range("J132").Value = 4057.1428
$ws.Range("K132").Value = 220620.9
$ws.Range("L132").Value = 12171.4284
$ws.Range("M132").Value = -218090.9
$ws.Range("N132").Value = -17231.4284
$ws.Range("H136").Value = 4985.143
$ws.Range("I136").Value = 3979.2
$ws.Range("J136").Value = 7500
$ws.Range("K136").Value = 11937.6
$ws.Range("L136").Value = 22500
$ws.Range("M136").Value = -9387.599999999999
$ws.Range("N136").Value = -27600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 74813
$ws.Range("J27").Value = 74813
$ws.Range("L27").Value = 74813
$ws.Range("N27").Value = -74951
$ws.Range("H54").Value = 23600
$ws.Range("I54").Value = 19750
$ws.Range("K54").Value = 19750
$ws.Range("M54").Value = -19230
$ws.Range("H58").Value = 40028
$ws.Range("I58").Value = 35042.5
$ws.Range("J58").Value = 49999
$ws.Range("K58").Value = 35042.5
$ws.Range("L58").Value = 49999
$ws.Range("M58").Value = -34734.5
$ws.Range("N58").Value = -50615
$ws.Range("H107").Value = 1121.1111
$ws.Range("J107").Value = 1229.6666
$ws.Range("L107").Value = 3688.9998
$ws.Range("N107").Value = -7528.9998
$ws.Range("H122").Value = 2875.9092
$ws.Range("I122").Value = 2859.4443
$ws.Range("K122").Value = 8578.332900000001
$ws.Range("M122").Value = -6128.332900000001
$ws.Range("H132").Value = 557786.7
$ws.Range("I132").Value = 760401.4
$ws.Range("J132").Value = 51250
$ws.Range("K132").Value = 2281204.2
$ws.Range("L132").Value = 153750
$ws.Range("M132").Value = -2278674.2
$ws.Range("N132").Value = -158810
